# SPFx Webinar deck - "CD/CD" -> "CI/CD" correction pass (2020-05-27).
#
# The slide titled "Continuous Integration Continuous Deployment (CD/CD)"
# (slide 8, sldId 1600) had a typo in the CI/CD acronym, and the adjoining
# SmartArt graphic still labelled its second stage "Continuous Development"
# instead of "Continuous Deployment". Both get corrected here.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

# --- Fix the SmartArt node text: "Continuous Development" -> "Continuous Deployment"
$graphicFrame = $s.Shapes.Item(1)
if ($graphicFrame.HasSmartArt) {
    $smartArt = $graphicFrame.SmartArt
    for ($i = 1; $i -le $smartArt.AllNodes.Count; $i++) {
        $node = $smartArt.AllNodes.Item($i)
        if ($node.TextFrame2.TextRange.Text -eq "Continuous Development") {
            $node.TextFrame2.TextRange.Text = "Continuous Deployment"
        }
    }
}

# --- Fix the slide title typo: "(CD/CD)" -> "(CI/CD)"
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $txt = $shp.TextFrame.TextRange.Text
        if ($txt -like "*CD/CD*") {
            $shp.TextFrame.TextRange.Text = $txt -replace "CD/CD", "CI/CD"
        }
    }
}
